$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update bracket_y (mm) value in B6 from 30 to 27
$ws.Range("B6").Value = 27

# Add new parameter rows 15-17
$ws.Range("A15").Value = "Rcaliper_offset (mm)"
$ws.Range("B15").Value = 12

$ws.Range("A16").Value = "Rcaliper_radius (mm)"
$ws.Range("B16").Value = 92

$ws.Range("A17").Value = "Rcaliper_distance (mm)"
$ws.Range("B17").Value = 84

# Update the selection to match the target state
$ws.Range("B20").Select()
